$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-5 from 2023-10-09 (45208) to 2023-10-13 (45212)
$ws.Range("C2").Value = 45212
$ws.Range("C3").Value = 45212
$ws.Range("C4").Value = 45212
$ws.Range("C5").Value = 45212
